$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Timp3"
$ws.Range("C2").Value = "Kdr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 150.258513
$ws.Range("H2").Value = 450.775539
$ws.Range("I2").Value = 0.5395416880146598
$ws.Range("J2").Value = 0.5395416880146598
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 161.7750676666667
$ws.Range("N2").Value = 485.325203
$ws.Range("O2").Value = 0.9790864123038654
$ws.Range("P2").Value = 0.9790864123038654
$ws.Range("Q2").Value = 24308.08110806771
$ws.Range("R2").Value = 218772.7299726094
$ws.Range("S2").Value = 0.5282579356066448
$ws.Range("T2").Value = 0.5282579356066448

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Timp3"
$ws.Range("C3").Value = "Kdr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 150.258513
$ws.Range("H3").Value = 450.775539
$ws.Range("I3").Value = 0.5395416880146598
$ws.Range("J3").Value = 0.5395416880146598
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.67894
$ws.Range("N3").Value = 2.03682
$ws.Range("O3").Value = 0.004109044356199978
$ws.Range("P3").Value = 0.004109044356199979
$ws.Range("Q3").Value = 102.01651481622
$ws.Range("R3").Value = 918.14863334598
$ws.Range("S3").Value = 0.002217000728071248
$ws.Range("T3").Value = 0.002217000728071248

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Timp3"
$ws.Range("C4").Value = "Kdr"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 150.258513
$ws.Range("H4").Value = 450.775539
$ws.Range("I4").Value = 0.5395416880146598
$ws.Range("J4").Value = 0.5395416880146598
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.763201333333333
$ws.Range("N4").Value = 5.289604
$ws.Range("O4").Value = 0.01067115280816804
$ws.Range("P4").Value = 0.01067115280816804
$ws.Range("Q4").Value = 264.936010466284
$ws.Range("R4").Value = 2384.424094196556
$ws.Range("S4").Value = 0.005757531799181362
$ws.Range("T4").Value = 0.005757531799181363

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Timp3"
$ws.Range("C5").Value = "Kdr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 150.258513
$ws.Range("H5").Value = 450.775539
$ws.Range("I5").Value = 0.5395416880146598
$ws.Range("J5").Value = 0.5395416880146598
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.013424
$ws.Range("N5").Value = 3.040272
$ws.Range("O5").Value = 0.006133390531766587
$ws.Range("P5").Value = 0.006133390531766588
$ws.Range("Q5").Value = 152.275583278512
$ws.Range("R5").Value = 1370.480249506608
$ws.Range("S5").Value = 0.003309219880762477
$ws.Range("T5").Value = 0.003309219880762477

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Timp3"
$ws.Range("C6").Value = "Kdr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 36.46294533333333
$ws.Range("H6").Value = 109.388836
$ws.Range("I6").Value = 0.1309295472339256
$ws.Range("J6").Value = 0.1309295472339256
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 161.7750676666667
$ws.Range("N6").Value = 485.325203
$ws.Range("O6").Value = 0.9790864123038654
$ws.Range("P6").Value = 0.9790864123038654
$ws.Range("Q6").Value = 5898.795448625967
$ws.Range("R6").Value = 53089.1590376337
$ws.Range("S6").Value = 0.1281913406658337
$ws.Range("T6").Value = 0.1281913406658337

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Timp3"
$ws.Range("C7").Value = "Kdr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 36.46294533333333
$ws.Range("H7").Value = 109.388836
$ws.Range("I7").Value = 0.1309295472339256
$ws.Range("J7").Value = 0.1309295472339256
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.67894
$ws.Range("N7").Value = 2.03682
$ws.Range("O7").Value = 0.004109044356199978
$ws.Range("P7").Value = 0.004109044356199979
$ws.Range("Q7").Value = 24.75615210461333
$ws.Range("R7").Value = 222.80536894152
$ws.Range("S7").Value = 0.0005379953171213807
$ws.Range("T7").Value = 0.0005379953171213808

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Timp3"
$ws.Range("C8").Value = "Kdr"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 36.46294533333333
$ws.Range("H8").Value = 109.388836
$ws.Range("I8").Value = 0.1309295472339256
$ws.Range("J8").Value = 0.1309295472339256
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.763201333333333
$ws.Range("N8").Value = 5.289604
$ws.Range("O8").Value = 0.01067115280816804
$ws.Range("P8").Value = 0.01067115280816804
$ws.Range("Q8").Value = 64.29151382899377
$ws.Range("R8").Value = 578.6236244609439
$ws.Range("S8").Value = 0.001397169205637476
$ws.Range("T8").Value = 0.001397169205637476

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Timp3"
$ws.Range("C9").Value = "Kdr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 36.46294533333333
$ws.Range("H9").Value = 109.388836
$ws.Range("I9").Value = 0.1309295472339256
$ws.Range("J9").Value = 0.1309295472339256
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.013424
$ws.Range("N9").Value = 3.040272
$ws.Range("O9").Value = 0.006133390531766587
$ws.Range("P9").Value = 0.006133390531766588
$ws.Range("Q9").Value = 36.95242391148799
$ws.Range("R9").Value = 332.571815203392
$ws.Range("S9").Value = 0.0008030420453330457
$ws.Range("T9").Value = 0.0008030420453330458

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Timp3"
$ws.Range("C10").Value = "Kdr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2335036666666667
$ws.Range("H10").Value = 0.700511
$ws.Range("I10").Value = 0.0008384547401380566
$ws.Range("J10").Value = 0.0008384547401380566
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 161.7750676666667
$ws.Range("N10").Value = 485.325203
$ws.Range("O10").Value = 0.9790864123038654
$ws.Range("P10").Value = 0.9790864123038654
$ws.Range("Q10").Value = 37.77507147541478
$ws.Range("R10").Value = 339.975643278733
$ws.Range("S10").Value = 0.0008209196434009397
$ws.Range("T10").Value = 0.0008209196434009397

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Timp3"
$ws.Range("C11").Value = "Kdr"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2335036666666667
$ws.Range("H11").Value = 0.700511
$ws.Range("I11").Value = 0.0008384547401380566
$ws.Range("J11").Value = 0.0008384547401380566
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.67894
$ws.Range("N11").Value = 2.03682
$ws.Range("O11").Value = 0.004109044356199978
$ws.Range("P11").Value = 0.004109044356199979
$ws.Range("Q11").Value = 0.1585349794466667
$ws.Range("R11").Value = 1.42681481502
$ws.Range("S11").Value = 0.000003445247717893401
$ws.Range("T11").Value = 0.000003445247717893402

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Timp3"
$ws.Range("C12").Value = "Kdr"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2335036666666667
$ws.Range("H12").Value = 0.700511
$ws.Range("I12").Value = 0.0008384547401380566
$ws.Range("J12").Value = 0.0008384547401380566
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.763201333333333
$ws.Range("N12").Value = 5.289604
$ws.Range("O12").Value = 0.01067115280816804
$ws.Range("P12").Value = 0.01067115280816804
$ws.Range("Q12").Value = 0.4117139764048888
$ws.Range("R12").Value = 3.705425787644
$ws.Range("S12").Value = 0.000008947278654746027
$ws.Range("T12").Value = 0.000008947278654746029

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Timp3"
$ws.Range("C13").Value = "Kdr"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2335036666666667
$ws.Range("H13").Value = 0.700511
$ws.Range("I13").Value = 0.0008384547401380566
$ws.Range("J13").Value = 0.0008384547401380566
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.013424
$ws.Range("N13").Value = 3.040272
$ws.Range("O13").Value = 0.006133390531766587
$ws.Range("P13").Value = 0.006133390531766588
$ws.Range("Q13").Value = 0.236638219888
$ws.Range("R13").Value = 2.129743978992
$ws.Range("S13").Value = 0.000005142570364477571
$ws.Range("T13").Value = 0.000005142570364477572

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Timp3"
$ws.Range("C14").Value = "Kdr"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 91.53790766666667
$ws.Range("H14").Value = 274.613723
$ws.Range("I14").Value = 0.3286903100112765
$ws.Range("J14").Value = 0.3286903100112764
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 161.7750676666667
$ws.Range("N14").Value = 485.325203
$ws.Range("O14").Value = 0.9790864123038654
$ws.Range("P14").Value = 0.9790864123038654
$ws.Range("Q14").Value = 14808.55120684009
$ws.Range("R14").Value = 133276.9608615608
$ws.Range("S14").Value = 0.321816216387986
$ws.Range("T14").Value = 0.3218162163879859

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Timp3"
$ws.Range("C15").Value = "Kdr"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 91.53790766666667
$ws.Range("H15").Value = 274.613723
$ws.Range("I15").Value = 0.3286903100112765
$ws.Range("J15").Value = 0.3286903100112764
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.67894
$ws.Range("N15").Value = 2.03682
$ws.Range("O15").Value = 0.004109044356199978
$ws.Range("P15").Value = 0.004109044356199979
$ws.Range("Q15").Value = 62.14874703120667
$ws.Range("R15").Value = 559.33872328086
$ws.Range("S15").Value = 0.001350603063289457
$ws.Range("T15").Value = 0.001350603063289457

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Timp3"
$ws.Range("C16").Value = "Kdr"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 91.53790766666667
$ws.Range("H16").Value = 274.613723
$ws.Range("I16").Value = 0.3286903100112765
$ws.Range("J16").Value = 0.3286903100112764
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.763201333333333
$ws.Range("N16").Value = 5.289604
$ws.Range("O16").Value = 0.01067115280816804
$ws.Range("P16").Value = 0.01067115280816804
$ws.Range("Q16").Value = 161.3997608484102
$ws.Range("R16").Value = 1452.597847635692
$ws.Range("S16").Value = 0.003507504524694457
$ws.Range("T16").Value = 0.003507504524694457

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Timp3"
$ws.Range("C17").Value = "Kdr"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 91.53790766666667
$ws.Range("H17").Value = 274.613723
$ws.Range("I17").Value = 0.3286903100112765
$ws.Range("J17").Value = 0.3286903100112764
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.013424
$ws.Range("N17").Value = 3.040272
$ws.Range("O17").Value = 0.006133390531766587
$ws.Range("P17").Value = 0.006133390531766588
$ws.Range("Q17").Value = 92.76671253918398
$ws.Range("R17").Value = 834.9004128526559
$ws.Range("S17").Value = 0.002015986035306588
$ws.Range("T17").Value = 0.002015986035306588
